$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gained one new weekly data row. It was inserted as row 16,
# pushing the previous rows 16-73 down to 17-74 (their contents/formatting
# unchanged - Insert() carries the row-above formatting onto the new row,
# which keeps the Fecha column's date number format on D16).
$ws.Rows("16:16").Insert()

$ws.Range("A16").Value2 = 4
$ws.Range("B16").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C16").Value2 = "Los Lagos"
$ws.Range("D16").Value2 = 44687
$ws.Range("E16").Value2 = 10
$ws.Range("F16").Value2 = 100112031
$ws.Range("G16").Value2 = "Poroto verde"
$ws.Range("H16").Value2 = "Magnum"
$ws.Range("I16").Value2 = "Primera"
$ws.Range("J16").Value2 = 20
$ws.Range("K16").Value2 = 35000
$ws.Range("L16").Value2 = 35000
$ws.Range("M16").Value2 = 35000
$ws.Range("N16").Value2 = "$/saco 25 kilos"
$ws.Range("O16").Value2 = "Región Metropolitana"
$ws.Range("P16").Value2 = 1400
$ws.Range("Q16").Value2 = 25
$ws.Range("R16").Value2 = "Hortaliza"
